# Add a new data row (row 56) to the worksheet, mirroring the existing
# rows' structure (Hortaliza / Vega Modelo de Temuco - Rabanito dataset).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A56").Value = 10
$ws.Range("B56").Value = "Vega Modelo de Temuco"
$ws.Range("C56").Value = "La Araucanía"
$ws.Range("D56").Value = 44595
$ws.Range("D56").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E56").Value = 9
$ws.Range("F56").Value = 300000001
$ws.Range("G56").Value = "Rabanito"
$ws.Range("H56").Value = "Sin especificar"
$ws.Range("I56").Value = "Primera"
$ws.Range("J56").Value = 30
$ws.Range("K56").Value = 6000
$ws.Range("L56").Value = 6000
$ws.Range("M56").Value = 6000
$ws.Range("N56").Value = "$/docena de paquetes"
$ws.Range("O56").Value = "Provincia de Cautín"
$ws.Range("P56").Value = 500
$ws.Range("Q56").Value = 12
$ws.Range("R56").Value = "Hortaliza"
